$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.572.93'
$ws.Range("E2").Value = '  +0.90%  '
$ws.Range("D3").Value = '1.873.01'
$ws.Range("E3").Value = '  +0.12%  '
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").Value = "'247.20"
$ws.Range("E5").Value = '  +0.81%  '
$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = '  +0.03%  '
$ws.Range("D7").Value = "'0.4738"
$ws.Range("E7").Value = '  +0.28%  '
$ws.Range("D8").Value = "'0.2908"
$ws.Range("E8").Value = '  +1.09%  '
$ws.Range("D9").Value = "'0.06477"
$ws.Range("E9").Value = '  +0.03%  '
$ws.Range("D10").Value = "'22.08"
$ws.Range("E10").Value = '  +4.63%  '
$ws.Range("D11").Value = "'0.07712"
$ws.Range("E11").Value = '  -0.61%  '
$ws.Range("D12").Value = "'0.7422"
$ws.Range("E12").Value = '  +4.25%  '
$ws.Range("D13").Value = "'96.54"
$ws.Range("E13").Value = '  +1.49%  '
$ws.Range("D14").Value = '1.870.98'
$ws.Range("E14").Value = '  -0.17%  '
$ws.Range("D15").Value = "'5.154"
$ws.Range("E15").Value = '  +1.00%  '
$ws.Range("D16").Value = "'273.50"
$ws.Range("E16").Value = '  -0.89%  '
$ws.Range("D17").Value = '30.567.02'
$ws.Range("E17").Value = '  +0.89%  '
$ws.Range("D18").Value = "'13.36"
$ws.Range("E18").Value = '  +0.19%  '
$ws.Range("D19").Value = "'1.000"
$ws.Range("E19").Value = '  +0.00%  '
$ws.Range("D20").Value = "'0.000007503"
$ws.Range("E20").Value = '  -0.92%  '
$ws.Range("D21").Value = '2.114.50'
$ws.Range("E21").Value = '  -0.26%  '
$ws.Range("E22").Value = '  -0.02%  '
$ws.Range("D23").Value = "'5.260"
$ws.Range("E23").Value = '  +0.63%  '
$ws.Range("D24").Value = "'6.191"
$ws.Range("E24").Value = '  +0.63%  '
$ws.Range("D25").Value = "'9.238"
$ws.Range("D26").Value = "'163.34"
$ws.Range("D27").Value = "'18.80"
$ws.Range("E27").Value = '  -0.52%  '
$ws.Range("D28").Value = "'1.920"
$ws.Range("E28").Value = '  +0.70%  '
$ws.Range("D29").Value = "'0.09991"
$ws.Range("E29").Value = '  +1.26%  '
$ws.Range("D30").Value = "'1.347"
$ws.Range("E30").Value = '  -2.22%  '
$ws.Range("D31").Value = "'1.506"
$ws.Range("E31").Value = '  -0.80%  '
$ws.Range("D32").Value = "'4.297"
$ws.Range("E32").Value = '  +0.84%  '
$ws.Range("D33").Value = "'4.112"
$ws.Range("E33").Value = '  +2.01%  '
$ws.Range("D34").Value = "'0.04789"
$ws.Range("E35").Value = '  -0.12%  '
$ws.Range("D36").Value = "'0.6976"
$ws.Range("E36").Value = '  +0.48%  '
$ws.Range("D37").Value = "'0.9998"
$ws.Range("E37").Value = '  +0.05%  '
$ws.Range("D38").Value = "'2.717"
$ws.Range("E38").Value = '  +0.04%  '
$ws.Range("D39").Value = "'0.01850"
$ws.Range("E39").Value = '  -0.02%  '
$ws.Range("E40").Value = '  +0.18%  '
$ws.Range("D41").Value = "'6.201"
$ws.Range("E41").Value = '  -1.58%  '
$ws.Range("D42").Value = "'73.46"
$ws.Range("E42").Value = '  +4.06%  '
$ws.Range("D43").Value = "'1.969"
$ws.Range("E43").Value = '  +2.97%  '
$ws.Range("D44").Value = "'0.4182"
$ws.Range("E44").Value = '  +1.72%  '
$ws.Range("E45").Value = '  +0.06%  '
$ws.Range("D46").Value = "'0.8336"
$ws.Range("E46").Value = '  -1.10%  '
$ws.Range("D47").Value = "'102.55"
$ws.Range("E47").Value = '  +0.65%  '
$ws.Range("D48").Value = "'9.290"
$ws.Range("E48").Value = '  +0.46%  '
$ws.Range("D49").Value = "'35.38"
$ws.Range("E49").Value = '  +0.26%  '
$ws.Range("D50").Value = "'927.79"
$ws.Range("E50").Value = '  +0.86%  '
$ws.Range("D51").Value = "'6.964"
$ws.Range("E51").Value = '  -1.94%  '
